# BasicApp: update business data xlsx
# - Swap the "Julien Gonzalez" / "Renaud Joly" rows on the BasicUser sheet
# - Rename Chloé's surname from "Baffert Bui-Van" to "Bui-Van" (and derived
#   username / email)
# - Add a new BasicUser row for Roland Foucher
# - Refresh the mailto hyperlinks on column E to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BasicUser")

# --- 1. Swap row 4 (Julien Gonzalez) and row 5 (Renaud Joly) ---------------
$row4 = @($ws.Cells.Item(4,2).Value2, $ws.Cells.Item(4,3).Value2, $ws.Cells.Item(4,4).Value2, $ws.Cells.Item(4,5).Value2)
$row5 = @($ws.Cells.Item(5,2).Value2, $ws.Cells.Item(5,3).Value2, $ws.Cells.Item(5,4).Value2, $ws.Cells.Item(5,5).Value2)

$ws.Cells.Item(4,2).Value2 = $row5[0]
$ws.Cells.Item(4,3).Value2 = $row5[1]
$ws.Cells.Item(4,4).Value2 = $row5[2]
$ws.Cells.Item(4,5).Value2 = $row5[3]

$ws.Cells.Item(5,2).Value2 = $row4[0]
$ws.Cells.Item(5,3).Value2 = $row4[1]
$ws.Cells.Item(5,4).Value2 = $row4[2]
$ws.Cells.Item(5,5).Value2 = $row4[3]

# --- 2. Chloé: "Baffert Bui-Van" -> "Bui-Van" -------------------------------
$ws.Cells.Item(12,3).Value2 = "Bui-Van"
$ws.Cells.Item(12,4).Value2 = "cbuivan"
$ws.Cells.Item(12,5).Value2 = "chloe.buivan@kobalt.fr"

# --- 3. New row 15: Roland Foucher -----------------------------------------
# Copy row 14 first so the new row inherits the same number formats / styles,
# then overwrite the person-specific fields (id, firstName, lastName,
# username, email). password/enabled/groups/authorities/locale stay the same
# default values as every other active BasicUser row.
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).PasteSpecial()

$ws.Cells.Item(15,1).Value2 = 14
$ws.Cells.Item(15,2).Value2 = "Roland"
$ws.Cells.Item(15,3).Value2 = "Foucher"
$ws.Cells.Item(15,4).Value2 = "rfoucher"
$ws.Cells.Item(15,5).Value2 = "roland.foucher@kobalt.fr"

# --- 4. Rebuild the column-E mailto hyperlinks so they track the data ------
$ws.Cells.Item(2,5).Hyperlinks.Delete()

$emails = @(
  "laurent.almeras@kobalt.fr",
  "julien.benichou@kobalt.fr",
  "renaud.joly@kobalt.fr",
  "julien.gonzalez@kobalt.fr",
  "florian.lacreuse@kobalt.fr",
  "margot.piva@kobalt.fr",
  "mathieu.palley@kobalt.fr",
  "alexandre.wallois@kobalt.fr",
  "corinne.fagno@kobalt.fr",
  "anais.rouviere@kobalt.fr",
  "chloe.buivan@kobalt.fr",
  "aurelien.jolivet@kobalt.fr",
  "vincent.weber@kobalt.fr",
  "roland.foucher@kobalt.fr"
)

for ($i = 0; $i -lt $emails.Length; $i++) {
  $row = $i + 2
  $email = $emails[$i]
  $ws.Hyperlinks.Add($ws.Cells.Item($row,5), ("mailto:" + $email), "", "", $email)
}

# --- 5. Move the active selection (cosmetic, matches the saved view) -------
$ws.Range("L9").Select()
